$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update cell B11 so its (shared) string value becomes "1" instead of "R40".
# Force text storage so Excel keeps it as a shared string rather than a number.
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"
